$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$ws1 = $wb.Worksheets.Item("Metadata")
$ws1.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/encounter-class"
$ws1.Range("B3").Value = "8.0.0"
$ws1.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$ws1.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Encounter Class " ---
$ws3 = $wb.Worksheets.Item("Include from Encounter Class ")
$ws3.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/encounter-class"
